# gym_data.xlsx - Stray Dog scraper fix + workflow update
#
# 1) "4x4 Squat Racks" sheet: row 4 (Rogue RM-3 Monster Rack 2.0) had its
#    Image URL (E4) replaced by a scraped base64 data-URI image instead of
#    a remote link, so the hyperlink + hyperlink style on E4 go away too.
# 2) "Squat Stands" sheet: row 3 (Titan X-3 Series Tall Squat Stand) price
#    is now available.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "4x4 Squat Racks"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("4x4 Squat Racks")

# Bank the original hyperlink-cell look (blue/underline direct formatting)
# on a scratch cell so it can be re-applied after rebuilding the links -
# Hyperlinks.Add() stamps its own built-in "Hyperlink" style on whatever
# cell it targets, which would otherwise clobber the existing formatting
# of every cell except E4.
$ws1.Range("F4").Copy()
$ws1.Range("Z1").PasteSpecial(-4122)

# This engine's Hyperlinks.Delete() (whether called on a Range or the
# Worksheet) always clears every hyperlink on the sheet - it is not
# scoped to a single cell - so rebuild the surviving links from scratch.
# Re-adding them reassigns rId1, rId2, ... in call order, which is
# exactly the renumbering the target file shows once E4's link is gone.
$ws1.Hyperlinks.Delete()

$ws1.Hyperlinks.Add($ws1.Range("E2"), "https://oakclubmfg.com/cdn/shop/products/210123_stockracks_8144_1800x1800.jpg?v=1617242493")
$ws1.Hyperlinks.Add($ws1.Range("F2"), "https://oakclubmfg.com/collections/corporate-racks/products/the-corporate-rack")
$ws1.Hyperlinks.Add($ws1.Range("E3"), "https://titan.fitness/cdn/shop/files/401223_01.jpg?v=1722443777&width=1946")
$ws1.Hyperlinks.Add($ws1.Range("F3"), "https://titan.fitness/products/titan-series-power-rack-90-36?variant=47930285916437")
$ws1.Hyperlinks.Add($ws1.Range("F4"), "https://www.roguefitness.com/rm-3-bolt-together-monster-rack-2-0")
$ws1.Hyperlinks.Add($ws1.Range("E5"), "https://shop.straydogstrength.com/cdn/shop/files/2120-v2-FRAME.jpg?v=1739385447&width=1946")
$ws1.Hyperlinks.Add($ws1.Range("F5"), "https://shop.straydogstrength.com/products/alpha-half-rack")
$ws1.Hyperlinks.Add($ws1.Range("E6"), "https://cdn.shopify.com/s/files/1/2559/4942/products/XL_SingleRack_BlackTexture.210.jpg?v=1567697449")
$ws1.Hyperlinks.Add($ws1.Range("F6"), "https://www.sorinex.com/products/xl-half-rack?Attachment+Color=Black+Texture&Upgrades=None")

# Restore the banked formatting on every cell that still carries a link.
$ws1.Range("Z1").Copy()
$ws1.Range("E2").PasteSpecial(-4122)
$ws1.Range("F2").PasteSpecial(-4122)
$ws1.Range("E3").PasteSpecial(-4122)
$ws1.Range("F3").PasteSpecial(-4122)
$ws1.Range("F4").PasteSpecial(-4122)
$ws1.Range("E5").PasteSpecial(-4122)
$ws1.Range("F5").PasteSpecial(-4122)
$ws1.Range("E6").PasteSpecial(-4122)
$ws1.Range("F6").PasteSpecial(-4122)
$ws1.Range("Z1").Clear()

# E4 no longer links out to a remote image - the scraper now inlines the
# picture as a base64 data URI - so drop the hyperlink cell formatting
# back to Normal as well.
$ws1.Range("E4").Style = "Normal"
$ws1.Range("E4").Value = "data:image/jpeg;base64,/9j/4AAQSkZJRgABAQAAAQABAAD..."

# ---------------------------------------------------------------------
# Sheet 2: "Squat Stands"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Squat Stands")

# Price became available for the Titan X-3 Series Tall Squat Stand.
# Force text (matching the plain "$1,234.56" style text used by every
# other price cell in the sheet) instead of letting "$459.99" auto-
# convert into a currency number.
$ws2.Range("C3").NumberFormat = "@"
$ws2.Range("C3").Value = "$459.99"
$ws2.Range("C3").Style = "Normal"
